# Update workbook/sheet to reflect data refreshed through 2022-10-15
# (commit: "Add data for 2022-10-23")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab: "Through 2022-10-14" -> "Through 2022-10-15"
$ws.Name = "Through 2022-10-15"

# 2. Update the row label for October in column A (row 11)
$ws.Cells.Item(11, 1).Value = "October (through 10-15)"

# 3. Update the October row (row 11) figures for the columns that changed
$ws.Cells.Item(11, 3).Value = 26    # C11: 2016
$ws.Cells.Item(11, 4).Value = 29    # D11: 2017
$ws.Cells.Item(11, 5).Value = 38    # E11: 2018
$ws.Cells.Item(11, 7).Value = 75    # G11: 2020
$ws.Cells.Item(11, 8).Value = 95    # H11: 2021
$ws.Cells.Item(11, 9).Value = 48    # I11: 2022

# 4. Update the Total row (row 12) figures for the columns that changed
$ws.Cells.Item(12, 3).Value = 455    # C12: 2016
$ws.Cells.Item(12, 4).Value = 656    # D12: 2017
$ws.Cells.Item(12, 5).Value = 586    # E12: 2018
$ws.Cells.Item(12, 7).Value = 976    # G12: 2020
$ws.Cells.Item(12, 8).Value = 1342   # H12: 2021
$ws.Cells.Item(12, 9).Value = 1326   # I12: 2022
